$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = "glyphosate 1 week before maize planting; metalochlor, atrazine, and mesotrione at planting"
$ws.Range("J3").Value = "glyphosate 1 week before maize planting; metalochlor, atrazine, and mesotrione at planting"
$ws.Range("H6").Value = "glyphosate 1 week before maize planting; metalochlor, atrazine, and mesotrione at planting"
$ws.Range("J5").Value = "glyphosate 1 week before maize planting; metalochlor, atrazine, and mesotrione at planting"

$ws.Range("H7").Value = "glyphosate before planting; glyphosate and fluthiacet-methyl at planting"
$ws.Range("J7").Value = "glyphosate before planting; glyphosate and fluthiacet-methyl at planting"

$ws.Range("H8").Value = "glyphosate and acetochlor  before planting (April 15), atrazine, acetochlor at planting (May 14); acetochlor and glyphosate after planting (June 15)"
$ws.Range("J8").Value = "chlorimuron-ethyl, flumioxazin, pyroxasulfone, and glyphosate before planting, dicamba and acetochlor after planting"

$ws.Range("J8").Select()
